# "Updated readme and ordering info"
# BOM.xlsx: fill in the two missing MPN (column C) values that were left
# blank for the 100k/2512 resistor and the .2uF/0603 capacitor rows, and
# correct a transposed-digit typo in the Mini-Fit Jr. 8-pin connector's
# MPN (column C, row 11): 39293083 -> 39281083.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "100k ohm, 2512" resistor (R45, R46) -> MPN RMCF2512FT100K
$ws.Range("C6").Value = "RMCF2512FT100K"

# Row 7: ".2uF, 0603" capacitor (C1, C2) -> MPN 1206YC204KAT4A
$ws.Range("C7").Value = "1206YC204KAT4A"

# Row 11: Mini-Fit Jr. 8 Connections (J2) -> corrected MPN
$ws.Range("C11").Value = 39281083

# Leave the cursor on the last-edited cell, matching the saved selection.
$ws.Range("C7").Select()
